# Generate Report for Handoff
# Adds a new "fdef4bf7-95d4-49c4-b3b6-c00ce31e72ad" file row to each of the
# three report sheets (Overview, zh-cn, de-de), mirroring the existing
# "2794f42b-2333-43b0-bd4b-66da344643f6" row.

$wb = $excel.ActiveWorkbook

# Hyperlink font color FF6495ED (cornflower blue) expressed as a COM BGR long.
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------------
# Sheet "Overview" -> new row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "fdef4bf7-95d4-49c4-b3b6-c00ce31e72ad.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1f78f484215f7bf405e1cb422f2d6532926f2abb/e2e/fdef4bf7-95d4-49c4-b3b6-c00ce31e72ad.md", "", "", "fdef4bf7-95d4-49c4-b3b6-c00ce31e72ad.md")
$wsOverview.Range("A3").Font.Name = "Calibri"
$wsOverview.Range("A3").Font.Size = 11
$wsOverview.Range("A3").Font.Underline = 1
$wsOverview.Range("A3").Font.Color = $hyperlinkColor

$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

$wsOverview.Range("D3").Value = "2016-03-23 16:42:05"
$wsOverview.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> new row 3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A3").Value = "fdef4bf7-95d4-49c4-b3b6-c00ce31e72ad.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1f78f484215f7bf405e1cb422f2d6532926f2abb/e2e/fdef4bf7-95d4-49c4-b3b6-c00ce31e72ad.md", "", "", "fdef4bf7-95d4-49c4-b3b6-c00ce31e72ad.md")
$wsZhCn.Range("A3").Font.Name = "Calibri"
$wsZhCn.Range("A3").Font.Size = 11
$wsZhCn.Range("A3").Font.Underline = 1
$wsZhCn.Range("A3").Font.Color = $hyperlinkColor

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"

$wsZhCn.Range("D3").Value = "fdef4bf7-95d4-49c4-b3b6-c00ce31e72ad.20ec26b1102e7ccca83cf138998d608526ea8170.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6f5b94484de87ada955d2ec2bfdcbb3d8504402b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/fdef4bf7-95d4-49c4-b3b6-c00ce31e72ad.20ec26b1102e7ccca83cf138998d608526ea8170.zh-cn.xlf", "", "", "fdef4bf7-95d4-49c4-b3b6-c00ce31e72ad.20ec26b1102e7ccca83cf138998d608526ea8170.zh-cn.xlf")
$wsZhCn.Range("D3").Font.Name = "Calibri"
$wsZhCn.Range("D3").Font.Size = 11
$wsZhCn.Range("D3").Font.Underline = 1
$wsZhCn.Range("D3").Font.Color = $hyperlinkColor

$wsZhCn.Range("E3").Value = "2016-03-23 16:41:59"
$wsZhCn.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("J3").Value = "Include"

# ---------------------------------------------------------------------------
# Sheet "de-de" -> new row 3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A3").Value = "fdef4bf7-95d4-49c4-b3b6-c00ce31e72ad.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1f78f484215f7bf405e1cb422f2d6532926f2abb/e2e/fdef4bf7-95d4-49c4-b3b6-c00ce31e72ad.md", "", "", "fdef4bf7-95d4-49c4-b3b6-c00ce31e72ad.md")
$wsDeDe.Range("A3").Font.Name = "Calibri"
$wsDeDe.Range("A3").Font.Size = 11
$wsDeDe.Range("A3").Font.Underline = 1
$wsDeDe.Range("A3").Font.Color = $hyperlinkColor

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"

$wsDeDe.Range("D3").Value = "fdef4bf7-95d4-49c4-b3b6-c00ce31e72ad.20ec26b1102e7ccca83cf138998d608526ea8170.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8429535a6e371dbf7635586e20a98d46ffa7f043/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/fdef4bf7-95d4-49c4-b3b6-c00ce31e72ad.20ec26b1102e7ccca83cf138998d608526ea8170.de-de.xlf", "", "", "fdef4bf7-95d4-49c4-b3b6-c00ce31e72ad.20ec26b1102e7ccca83cf138998d608526ea8170.de-de.xlf")
$wsDeDe.Range("D3").Font.Name = "Calibri"
$wsDeDe.Range("D3").Font.Size = 11
$wsDeDe.Range("D3").Font.Underline = 1
$wsDeDe.Range("D3").Font.Color = $hyperlinkColor

$wsDeDe.Range("E3").Value = "2016-03-23 16:42:05"
$wsDeDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("J3").Value = "Include"

Write-Output "Handoff report rows added."
